$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.503401875495911
$ws.Range("B1").Value = 2.100000858306885
$ws.Range("C1").Value = 2.489894866943359
$ws.Range("D1").Value = 2.972882032394409
$ws.Range("E1").Value = 2.442269563674927
